$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: change F3 first so the new "1.5x7" shared string is created
#     before any of the other brand-new strings below. ---
$ws.Range("F3").Value = "1.5x7"

# --- Duplicate row 4 (values + formatting) into a brand-new row 5, and set
#     its F column next so "-2.19x13.49" becomes the 2nd newly created
#     shared string. ---
$ws.Rows("4:4").Copy()
$ws.Rows("5:5").Insert(-4121)
$ws.Rows("5:5").RowHeight = 15.75
$ws.Range("F5").Value = "-2.19x13.49"

# F5 gets its own distinct look, different from the rest of the row
$ws.Range("F5").HorizontalAlignment = -4152
$ws.Range("F5").VerticalAlignment = -4160

# --- Row 2 ---
$ws.Range("A2").Value = ">=500"
$ws.Range("B2").Value = "<2"
$ws.Range("C2").Value = "<2"

# --- Row 3 (remaining cells) ---
$ws.Range("A3").Value = ">=500"
$ws.Range("B3").Value = "<2"
$ws.Range("C3").Value = "<2"
$ws.Range("E3").Value = ">=300"

# --- Row 4 ---
$ws.Range("A4").Value = "<500"
$ws.Range("B4").Value = ">=2"

# --- Row 5 (remaining cells) ---
$ws.Range("C5").Value = ">=2"

# --- Restore the saved selection/active cell ---
$ws.Range("E3").Select()
